{"js": "// Remove the trailing \"Ver no Jupiter...\" / \"\u00a9 2020 ...\" footer block\n// (and the blank paragraph right before it) that used to follow the\n// \"LOT2045: Biologia (Requisito fraco)\" paragraph at the end of the body.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the anchor paragraph (\"LOT2045: ...\") and the two distinctly\n// worded paragraphs that must disappear, then delete everything from the\n// blank paragraph right after the anchor through the copyright paragraph\n// (inclusive).\nlet anchorIndex = -1;\nlet jupiterIndex = -1;\nlet copyrightIndex = -1;\n\nfor (let i = 0; i < items.length; i++) {\n  const text = items[i].text;\n  if (text.indexOf(\"LOT2045\") !== -1) {\n    anchorIndex = i;\n  }\n  if (text.indexOf(\"Ver no Jupiter\") !== -1) {\n    jupiterIndex = i;\n  }\n  if (text.indexOf(\"Powered by Jekyll\") !== -1) {\n    copyrightIndex = i;\n  }\n}\n\nif (jupiterIndex === -1 || copyrightIndex === -1) {\n  throw new Error(\"Could not locate the footer paragraphs to remove.\");\n}\n\n// The blank paragraph that separates \"LOT2045...\" from \"Ver no Jupiter...\"\n// is also removed (it sat directly above the footer block).\nlet startIndex = jupiterIndex - 1;\nif (anchorIndex !== -1 && startIndex <= anchorIndex) {\n  startIndex = anchorIndex + 1;\n}\n\nconst toDelete = [];\nfor (let i = startIndex; i <= copyrightIndex; i++) {\n  toDelete.push(items[i]);\n}\n\nfor (const paragraph of toDelete) {\n  paragraph.delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the trailing \"Ver no Jupiter...\" / \"(c) 2020 ...\" footer block\n# (and the blank paragraph right before it) that used to follow the\n# \"LOT2045: Biologia (Requisito fraco)\" paragraph at the end of the body.\n\n$d = $word.ActiveDocument\n\n$anchorIndex = -1\n$jupiterIndex = -1\n$copyrightIndex = -1\n$count = $d.Paragraphs.Count\n\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs($i).Range.Text\n    if ($t -like \"*LOT2045*\") { $anchorIndex = $i }\n    if ($t -like \"*Ver no Jupiter*\") { $jupiterIndex = $i }\n    if ($t -like \"*Powered by Jekyll*\") { $copyrightIndex = $i }\n}\n\nif ($jupiterIndex -eq -1 -or $copyrightIndex -eq -1) {\n    throw \"Could not locate the footer paragraphs to remove.\"\n}\n\n# The blank paragraph that separates \"LOT2045...\" from \"Ver no Jupiter...\"\n# is also removed (it sat directly above the footer block).\n$startIndex = $jupiterIndex - 1\nif ($anchorIndex -ne -1 -and $startIndex -le $anchorIndex) {\n    $startIndex = $anchorIndex + 1\n}\n\n$rangeStart = $d.Paragraphs($startIndex).Range.Start\n$rangeEnd = $d.Paragraphs($copyrightIndex).Range.End\n\n$deleteRange = $d.Range($rangeStart, $rangeEnd)\n$deleteRange.Delete()\n"}
